$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns Q:S entirely (dimension shrinks from A1:S8 to A1:P8)
$ws.Range("Q1:S8").EntireColumn.Delete()

# Row 2 - Task
$ws.Range("B2:P2").Value = "Table_VQA"

# Row 3 - Input Type
$ws.Range("B3").Value = "text"
$ws.Range("C3").Value = "text"
$ws.Range("D3").Value = "text"
$ws.Range("E3").Value = "text"
$ws.Range("F3").Value = "text"
$ws.Range("G3").Value = "image"
$ws.Range("H3").Value = "image"
$ws.Range("I3").Value = "image"
$ws.Range("J3").Value = "image"
$ws.Range("K3").Value = "image"
$ws.Range("L3").Value = "hybrid"
$ws.Range("M3").Value = "hybrid"
$ws.Range("N3").Value = "hybrid"
$ws.Range("O3").Value = "hybrid"
$ws.Range("P3").Value = "hybrid"

# Row 4 - Model Type
$ws.Range("B4").Value = "gpt-4o"
$ws.Range("C4").Value = "claude"
$ws.Range("D4").Value = "google"
$ws.Range("E4").Value = "qwen-25"
$ws.Range("F4").Value = "gemini-1.5-pro"
$ws.Range("G4").Value = "gpt-4o"
$ws.Range("H4").Value = "claude"
$ws.Range("I4").Value = "google"
$ws.Range("J4").Value = "qwen-25"
$ws.Range("K4").Value = "gemini-1.5-pro"
$ws.Range("L4").Value = "gpt-4o"
$ws.Range("M4").Value = "claude"
$ws.Range("N4").Value = "google"
$ws.Range("O4").Value = "qwen-25"
$ws.Range("P4").Value = "gemini-1.5-pro"

# Row 5 - Exact Match Mean
$ws.Range("B5").Value = 0.24
$ws.Range("C5").Value = 0.224
$ws.Range("D5").Value = 0.112
$ws.Range("E5").Value = 0.232
$ws.Range("F5").Value = 0.124
$ws.Range("G5").Value = 0.452
$ws.Range("H5").Value = 0.444
$ws.Range("I5").Value = 0.2409638554216867
$ws.Range("J5").Value = 0.004016064257028112
$ws.Range("K5").Value = 0.104
$ws.Range("L5").Value = 0.476
$ws.Range("M5").Value = 0.504
$ws.Range("N5").Value = 0.224
$ws.Range("O5").Value = 0.2530120481927711
$ws.Range("P5").Value = 0.176

# Row 6 - F1-Score Mean
$ws.Range("B6").Value = 0.4122325814536341
$ws.Range("C6").Value = 0.4548349206349206
$ws.Range("D6").Value = 0.3352632622279681
$ws.Range("E6").Value = 0.4044222222222222
$ws.Range("F6").Value = 0.3475746031746031
$ws.Range("G6").Value = 0.7835238095238094
$ws.Range("H6").Value = 0.8113142857142857
$ws.Range("I6").Value = 0.5797956007214092
$ws.Range("J6").Value = 0.06601644673933831
$ws.Range("K6").Value = 0.3539272727272727
$ws.Range("L6").Value = 0.8014545454545455
$ws.Range("M6").Value = 0.7983437229437229
$ws.Range("N6").Value = 0.5556857142857143
$ws.Range("O6").Value = 0.4312009944540065
$ws.Range("P6").Value = 0.5326666666666666

# Row 7 - METEOR Mean
$ws.Range("B7").Value = 0.08043396517781895
$ws.Range("C7").Value = 0.09751345912408725
$ws.Range("D7").Value = 0.07924119985200585
$ws.Range("E7").Value = 0.08500407969758361
$ws.Range("F7").Value = 0.1026602918901738
$ws.Range("G7").Value = 0.1483364357530102
$ws.Range("H7").Value = 0.1666057752180688
$ws.Range("I7").Value = 0.117876441862915
$ws.Range("J7").Value = 0.03980397267895593
$ws.Range("K7").Value = 0.1017992816128101
$ws.Range("L7").Value = 0.09963697991820604
$ws.Range("M7").Value = 0.1240272302853829
$ws.Range("N7").Value = 0.08884722823956065
$ws.Range("O7").Value = 0.06591144205267524
$ws.Range("P7").Value = 0.09209420530228477

# Row 8 - label change + Bert Score Mean values
$ws.Range("A8").Value = "Bert Score Mean"
$ws.Range("B8").Value = 0.6631027999520301
$ws.Range("C8").Value = 0.6769690741896629
$ws.Range("D8").Value = 0.5380090215802192
$ws.Range("E8").Value = 0.666713973402977
$ws.Range("F8").Value = 0.5734538987874985
$ws.Range("G8").Value = 0.8502521013021469
$ws.Range("H8").Value = 0.8572390868663787
$ws.Range("I8").Value = 0.7117472038450969
$ws.Range("J8").Value = 0.6062006224470445
$ws.Range("K8").Value = 0.6983631743192673
$ws.Range("L8").Value = 0.8581778284311294
$ws.Range("M8").Value = 0.8664072604179383
$ws.Range("N8").Value = 0.6925916314125061
$ws.Range("O8").Value = 0.723248125259656
$ws.Range("P8").Value = 0.7167851884961128
